$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct the quiz answer grid: "-4" responses relaxed to "-3", and several
#     other answers nudged one step higher, trimming total error by a few points ---
$ws.Range("C3").Value = -3
$ws.Range("D3").Value = -1
$ws.Range("B4").Value = -2
$ws.Range("D4").Value = -3
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = -1
$ws.Range("C7").Value = -3
$ws.Range("D8").Value = -1
$ws.Range("D11").Value = -3
$ws.Range("E11").Value = -1
$ws.Range("D12").Value = -3
$ws.Range("F12").Value = -1
$ws.Range("D14").Value = -3
$ws.Range("F14").Value = -3
$ws.Range("C16").Value = -3
$ws.Range("D16").Value = -3
$ws.Range("G16").Value = -2
$ws.Range("C17").Value = -3
$ws.Range("D17").Value = -3
$ws.Range("D18").Value = -1
$ws.Range("C19").Value = -3
$ws.Range("D19").Value = -3
$ws.Range("C20").Value = -3
$ws.Range("D20").Value = -3
$ws.Range("D22").Value = -3
$ws.Range("F22").Value = -1
$ws.Range("D24").Value = -1
$ws.Range("D25").Value = -1
$ws.Range("D26").Value = -3
$ws.Range("D27").Value = -3
$ws.Range("D28").Value = -3
$ws.Range("E28").Value = -3
$ws.Range("C29").Value = -3
$ws.Range("D29").Value = -3
$ws.Range("D31").Value = -3
$ws.Range("C32").Value = -3
$ws.Range("D34").Value = -3
$ws.Range("F34").Value = -3
$ws.Range("D35").Value = -3
$ws.Range("D37").Value = -1
$ws.Range("D38").Value = -3
$ws.Range("F38").Value = -1
$ws.Range("D39").Value = -3
$ws.Range("G39").Value = -2
$ws.Range("D40").Value = -3
$ws.Range("D42").Value = -1
$ws.Range("D44").Value = -1
$ws.Range("D46").Value = -1
$ws.Range("D49").Value = -3
$ws.Range("D52").Value = -1
$ws.Range("D53").Value = -3
$ws.Range("E53").Value = -1
$ws.Range("A64").Value = -1
$ws.Range("D64").Value = -3
$ws.Range("E64").Value = -2
$ws.Range("A65").Value = -1
$ws.Range("B65").Value = -3
$ws.Range("D65").Value = -3
$ws.Range("D66").Value = -3
$ws.Range("C69").Value = -3
$ws.Range("D69").Value = -3
$ws.Range("E69").Value = -1
$ws.Range("D70").Value = -2
$ws.Range("E70").Value = -3
$ws.Range("C71").Value = -3
$ws.Range("A72").Value = -1
$ws.Range("C72").Value = -3
$ws.Range("D72").Value = -3
$ws.Range("E72").Value = -2
$ws.Range("G72").Value = -2
$ws.Range("C73").Value = -3
$ws.Range("D73").Value = -3
$ws.Range("H73").Value = -2
$ws.Range("C75").Value = -3
$ws.Range("A77").Value = 0
$ws.Range("C77").Value = -3
$ws.Range("D77").Value = -3
$ws.Range("D79").Value = -3
$ws.Range("C80").Value = -3
$ws.Range("E80").Value = -1
$ws.Range("D81").Value = -1

# --- Append one more respondent (row 82), then extend the soma/resp/duplicate
#     check formulas down to cover it, matching the pattern used through row 81 ---
$ws.Range("A82").Value = 0
$ws.Range("B82").Value = -1
$ws.Range("C82").Value = -3
$ws.Range("D82").Value = -2
$ws.Range("E82").Value = -3
$ws.Range("F82").Value = -1
$ws.Range("G82").Value = -3
$ws.Range("H82").Value = 0
$ws.Range("I82").Formula = "=26+SUM(A82:H82)"
$ws.Range("K82").Formula = '=IF(I82<7,"0, 0, 0, 0, 1", IF(I82<13,"0, 0, 0, 1, 0", IF(I82<19,"0, 0, 1, 0, 0",IF(I82<25,"0, 1, 0, 0, 0",IF(I82<27,"1, 0, 0, 0, 0","")))))'
$ws.Range("L82").Formula = '=IF(SUMPRODUCT(($A$2:$A$81=A82)*1, ($B$2:$B$81=B82)*1, ($C$2:$C$81=C82)*1,($D$2:$D$81=D82)*1,($E$2:$E$81=E82)*1,($F$2:$F$81=F82)*1,($G$2:$G$81=G82)*1,($H$2:$H$81=H82)*1 ) >1, "duplicates", "no duplicates")'

# --- Leave the view scrolled near the bottom of the table, selection on the last edited cell ---
$ws.Range("C80").Select()
$excel.ActiveWindow.ScrollRow = 70
